$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 522, shifting the
# existing rows 522:605 down to 523:606 (dimension grows to A1:R606).
$ws.Rows("522").Insert()

# Populate the freshly inserted row 522 with the new weekly record.
$ws.Range("A522").Value = 10
$ws.Range("B522").Value = "Vega Modelo de Temuco"
$ws.Range("C522").Value = "La Araucanía"
$ws.Range("D522").Value = 44984
$ws.Range("E522").Value = 9
$ws.Range("F522").Value = 100112023
$ws.Range("G522").Value = "Brócoli"
$ws.Range("H522").Value = "Sin especificar"
$ws.Range("I522").Value = "Primera"
$ws.Range("J522").Value = 1500
$ws.Range("K522").Value = 1300
$ws.Range("L522").Value = 1300
$ws.Range("M522").Value = 1300
$ws.Range("N522").Value = "$/unidad"
$ws.Range("O522").Value = "Región de O'Higgins"
$ws.Range("P522").Value = 1300
$ws.Range("Q522").Value = 1
$ws.Range("R522").Value = "Hortaliza"
